# Weekly update: a new daily price record was reported for
# "Terminal Hortofrutícola Agro Chillán" - Zanahoria.
# It is inserted as the new row 101, pushing the previously existing
# rows 101-160 down to 102-161 (dimension grows from R160 to R161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 101 (shifts rows 101:160 -> 102:161)
$ws.Rows(101).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A101").Value = 7
$ws.Range("B101").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C101").Value = "Ñuble"
$ws.Range("D101").Value = 44438
$ws.Range("E101").Value = 16
$ws.Range("F101").Value = 100114013
$ws.Range("G101").Value = "Zanahoria"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 160
$ws.Range("K101").Value = 5500
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = 5750
$ws.Range("N101").Value = "$/saco 20 kilos"
$ws.Range("O101").Value = "Provincia de Diguillín"
$ws.Range("P101").Value = 288
$ws.Range("Q101").Value = 20
$ws.Range("R101").Value = "Hortaliza"
